# Add a new weekly record for "Cebollín" (Terminal La Palmera de La Serena)
# by inserting a new row at row 202 and shifting existing rows 202:249 down to 203:250.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 202 (pushes old rows 202-249 down to 203-250)
$ws.Rows.Item(202).Insert()

# Populate the new row 202 with the new record's data
$ws.Cells.Item(202, 1).Value = 8
$ws.Cells.Item(202, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(202, 3).Value = "Coquimbo"
$ws.Cells.Item(202, 4).Value = 44855
$ws.Cells.Item(202, 5).Value = 4
$ws.Cells.Item(202, 6).Value = 100112037
$ws.Cells.Item(202, 7).Value = "Cebollín"
$ws.Cells.Item(202, 8).Value = "Sin especificar"
$ws.Cells.Item(202, 9).Value = "Primera"
$ws.Cells.Item(202, 10).Value = 1200
$ws.Cells.Item(202, 11).Value = 1400
$ws.Cells.Item(202, 12).Value = 1600
$ws.Cells.Item(202, 13).Value = 1500
$ws.Cells.Item(202, 14).Value = "`$/paquete 6 unidades"
$ws.Cells.Item(202, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(202, 16).Value = 250
$ws.Cells.Item(202, 17).Value = 6
$ws.Cells.Item(202, 18).Value = "Hortaliza"

# Match style used for the Date column (D) on the rest of the sheet
$ws.Cells.Item(202, 4).NumberFormat = $ws.Cells.Item(203, 4).NumberFormat
